$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 299
$ws.Range("I18").Value = 299
$ws.Range("K18").Value = 299
$ws.Range("M18").Value = -15

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2464.1667
$ws.Range("J40").Value = 2550.625
$ws.Range("L40").Value = 2550.625
$ws.Range("N40").Value = -2900.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2787.7778
$ws.Range("I51").Value = 2145.5
$ws.Range("J51").Value = 2971.2856
$ws.Range("K51").Value = 2145.5
$ws.Range("L51").Value = 2971.2856
$ws.Range("M51").Value = -1661.5
$ws.Range("N51").Value = -3939.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 225.66667
$ws.Range("I55").Value = 178.5
$ws.Range("J55").Value = 320
$ws.Range("K55").Value = 178.5
$ws.Range("L55").Value = 320
$ws.Range("M55").Value = 35.5
$ws.Range("N55").Value = -748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2864.7222
$ws.Range("I70").Value = 2001.4286
$ws.Range("J70").Value = 3414.0908
$ws.Range("K70").Value = 6004.2858
$ws.Range("L70").Value = 10242.2724
$ws.Range("M70").Value = -5734.2858
$ws.Range("N70").Value = -10782.2724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2864.7222
$ws.Range("I73").Value = 2001.4286
$ws.Range("J73").Value = 3414.0908
$ws.Range("K73").Value = 6004.2858
$ws.Range("L73").Value = 10242.2724
$ws.Range("M73").Value = -5068.2858
$ws.Range("N73").Value = -12114.2724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 228999.95
$ws.Range("I132").Value = 253760.23
$ws.Range("K132").Value = 761280.6900000001
$ws.Range("M132").Value = -758750.6900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3363.1177
$ws.Range("I88").Value = 2199.1667
$ws.Range("J88").Value = 3998
$ws.Range("K88").Value = 2199.1667
$ws.Range("L88").Value = 3998
$ws.Range("M88").Value = -1793.1667
$ws.Range("N88").Value = -4810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3363.1177
$ws.Range("I91").Value = 2199.1667
$ws.Range("J91").Value = 3998
$ws.Range("K91").Value = 2199.1667
$ws.Range("L91").Value = 3998
$ws.Range("M91").Value = -795.1667000000002
$ws.Range("N91").Value = -6806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8206.75
$ws.Range("I86").Value = 2557.2856
$ws.Range("K86").Value = 2557.2856
$ws.Range("M86").Value = -1434.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 8206.75
$ws.Range("I89").Value = 2557.2856
$ws.Range("K89").Value = 12786.428
$ws.Range("M89").Value = -7170.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 4939
$ws.Range("J8").Value = 6712.857
$ws.Range("L8").Value = 6712.857
$ws.Range("N8").Value = -6992.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 12500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 272.46667
$ws.Range("I107").Value = 127.75
$ws.Range("J107").Value = 325.0909
$ws.Range("K107").Value = 127.75
$ws.Range("L107").Value = 325.0909
$ws.Range("M107").Value = 1792.25
$ws.Range("N107").Value = -4165.0909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1805.7273
$ws.Range("I122").Value = 1086.5
$ws.Range("J122").Value = 3064.375
$ws.Range("K122").Value = 3259.5
$ws.Range("L122").Value = 9193.125
$ws.Range("M122").Value = -809.5
$ws.Range("N122").Value = -14093.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2214.7659
$ws.Range("I132").Value = 1756.0555
$ws.Range("J132").Value = 3716
$ws.Range("K132").Value = 5268.166499999999
$ws.Range("L132").Value = 11148
$ws.Range("M132").Value = -2738.166499999999
$ws.Range("N132").Value = -16208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1133.3334
$ws.Range("I80").Value = 800
$ws.Range("K80").Value = 2400
$ws.Range("M80").Value = -1464

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1133.3334
$ws.Range("I83").Value = 800
$ws.Range("K83").Value = 7200
$ws.Range("M83").Value = -2520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 23745
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 29993.334
$ws.Range("K57").Value = 5000
$ws.Range("L57").Value = 29993.334
$ws.Range("M57").Value = -4180
$ws.Range("N57").Value = -31633.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5809.125
$ws.Range("I70").Value = 5964.696
$ws.Range("J70").Value = 5411.5557
$ws.Range("K70").Value = 5964.696
$ws.Range("L70").Value = 5411.5557
$ws.Range("M70").Value = -5694.696
$ws.Range("N70").Value = -5951.5557

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5809.125
$ws.Range("I73").Value = 5964.696
$ws.Range("J73").Value = 5411.5557
$ws.Range("K73").Value = 5964.696
$ws.Range("L73").Value = 5411.5557
$ws.Range("M73").Value = -5028.696
$ws.Range("N73").Value = -7283.5557

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2656.5
$ws.Range("I80").Value = 2351.25
$ws.Range("J80").Value = 3063.5
$ws.Range("K80").Value = 2351.25
$ws.Range("L80").Value = 3063.5
$ws.Range("M80").Value = -1353.25
$ws.Range("N80").Value = -5059.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2656.5
$ws.Range("I83").Value = 2351.25
$ws.Range("J83").Value = 3063.5
$ws.Range("K83").Value = 11756.25
$ws.Range("L83").Value = 15317.5
$ws.Range("M83").Value = -6764.25
$ws.Range("N83").Value = -25301.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8919.076999999999
$ws.Range("I22").Value = 966.6667
$ws.Range("J22").Value = 11304.8
$ws.Range("K22").Value = 966.6667
$ws.Range("L22").Value = 11304.8
$ws.Range("M22").Value = -671.6667
$ws.Range("N22").Value = -11894.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8919.076999999999
$ws.Range("I27").Value = 966.6667
$ws.Range("J27").Value = 11304.8
$ws.Range("K27").Value = 966.6667
$ws.Range("L27").Value = 11304.8
$ws.Range("M27").Value = -859.6667
$ws.Range("N27").Value = -11518.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2167.389
$ws.Range("I68").Value = 2078.2222
$ws.Range("J68").Value = 2256.5557
$ws.Range("K68").Value = 2078.2222
$ws.Range("L68").Value = 2256.5557
$ws.Range("M68").Value = -1329.2222
$ws.Range("N68").Value = -3754.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2167.389
$ws.Range("I71").Value = 2078.2222
$ws.Range("J71").Value = 2256.5557
$ws.Range("K71").Value = 10391.111
$ws.Range("L71").Value = 11282.7785
$ws.Range("M71").Value = -6647.111000000001
$ws.Range("N71").Value = -18770.7785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25028200
$ws.Range("I62").Value = 33368766
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 33368766
$ws.Range("L62").Value = 6500
$ws.Range("M62").Value = -33368142
$ws.Range("N62").Value = -7748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 25028200
$ws.Range("I65").Value = 33368766
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 166843830
$ws.Range("L65").Value = 32500
$ws.Range("M65").Value = -166840710
$ws.Range("N65").Value = -38740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13160110
$ws.Range("I132").Value = 17243410
$ws.Range("J132").Value = 2811.889
$ws.Range("K132").Value = 51730230
$ws.Range("L132").Value = 8435.667000000001
$ws.Range("M132").Value = -51727700
$ws.Range("N132").Value = -13495.667
